$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5").Value = "{{item.Clave}}"
$ws.Range("C5").Value = "{{item.Nombre}}"
$ws.Range("D5").Value = "{{item.Registro}}"
$ws.Range("E5").Value = "{{item.Entrega}}"
$ws.Range("F5").Value = "{{item.NombreEstatus}}"

[void]$ws.Range("B33").Select()
